# edit.ps1 - applies the two rule edits described by the diff:
#   1. Split/rewrite the scoring-stop sentence: the trailing clause
#      "до завершения игры по желанию одного из игроков." becomes
#      "по достижению определенного количества очков у одного из игроков."
#   2. Merge/trim the "error" rules list:
#        - the "Если вводится..." rule is rewritten to read what used to
#          be the separate "Ввод несуществующего..." rule's text
#        - the old separate "Ввод несуществующего..." paragraph is removed
#          (its text now lives in the rewritten paragraph above)
#        - the "Ввод города без населения является ошибкой" paragraph is
#          removed entirely

$d = $word.ActiveDocument

# --- Change 1: rewrite the end of the "until N mistakes / ..." rule ---
$d.Content.Find.Execute(
    "до завершения игры по желанию одного из игроков.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "по достижению определенного количества очков у одного из игроков.",
    2
) | Out-Null

# --- Change 2: collapse the three "mistake" rules into one, drop the last ---

# Locate the three paragraphs by their distinctive text so this does not
# depend on absolute paragraph indices.
$targetRewrite = $null
$targetDupe = $null
$targetRemove = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Если вводится название города*") {
        $targetRewrite = $p
    } elseif ($t -like "Ввод несуществующего на данный момент города*") {
        $targetDupe = $p
    } elseif ($t -like "Ввод города без населения является ошибкой*") {
        $targetRemove = $p
    }
}

# Rewrite the first rule's text with the (soon to be de-duplicated) rule text.
$targetRewrite.Range.Text = "Ввод несуществующего на данный момент города (не неназванного переименованного, а несуществующего) является ошибкой."

# Remove the now-duplicate paragraph and the no-longer-wanted paragraph.
# Delete bottom-up so indices/ranges of the earlier paragraph stay valid.
$targetRemove.Range.Delete() | Out-Null
$targetDupe.Range.Delete() | Out-Null
